# Add a new "FRI headways and runtimes" tab, positioned right after the
# "WKDY headways and runtimes" tab and before "SAT headways and runtimes".
# The new tab is an exact duplicate of the WKDY tab's data/formatting
# (same header row, same 21 data rows, same footnote rows, same column
# widths and cell styles), matching how the original FRI fixture tab was
# produced (copied from the weekday sheet and renamed).

$wb = $excel.ActiveWorkbook

$wkdy = $wb.Worksheets.Item("WKDY headways and runtimes")

# Worksheet.Copy(Before, After) — place the copy immediately after WKDY,
# i.e. immediately before SAT. This duplicates all cell values, number
# formats, column widths, and styles from WKDY onto the new sheet.
$wkdy.Copy($null, $wkdy)

# The copy becomes the active sheet/tab, same as it would in real Excel.
$fri = $wb.ActiveSheet
$fri.Name = "FRI headways and runtimes"
